# TC07_CDS_Filter_PHSAccession-phs002250.xlsx
# "Filter - PHS Accession - Testcases"
#
# The underlying change is a re-save of the workbook from Excel: the wrapped
# query cells in column B/C reflow (rows 2-4 grow slightly taller) and a
# couple of the sheet's auto-fit columns re-measure a little wider. There is
# no data/content change - just the recomputed row heights / column widths
# that Excel stamps into the sheet on save. Re-apply the same measurements
# here.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Rows 2-4 hold the wrapped Cypher query text in B/C; their auto height grows
# slightly (171.6 -> 173.25, 218.4 -> 220.5).
$ws.Rows.Item(2).RowHeight = 173.25
$ws.Rows.Item(3).RowHeight = 173.25
$ws.Rows.Item(4).RowHeight = 220.5

# Column A is best-fit on the tab-name text; its measured width grows
# (12.33203125 -> ~14.7109375).
$ws.Columns.Item(1).ColumnWidth = 13.76

# Columns B/C are the (manually sized) query columns; they widen a hair
# (75.6640625 -> ~75.7109375).
$ws.Columns.Item(2).ColumnWidth = 74.76
$ws.Columns.Item(3).ColumnWidth = 74.76
